$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 4.759400000000001
$ws.Range("A9").Value = -21.92250000000002
$ws.Range("B12").Value = 5.270699999999998
$ws.Range("D13").Value = -8.516300000000001
$ws.Range("C15").Value = -13.37029999999999
$ws.Range("D16").Value = -8.845300000000003
$ws.Range("A18").Value = -22.15670000000001
$ws.Range("A20").Value = -20.05669999999998
$ws.Range("D20").Value = -7.144999999999997
$ws.Range("D24").Value = -7.4846
$ws.Range("B26").Value = 4.077500000000003
$ws.Range("A27").Value = -22.1683
$ws.Range("B27").Value = 5.486500000000003
$ws.Range("E27").Value = 16.75729999999999
$ws.Range("B29").Value = 5.069999999999995
$ws.Range("E29").Value = 17.17150000000002
$ws.Range("B37").Value = 9.877600000000003
$ws.Range("B38").Value = 4.897100000000002
$ws.Range("C38").Value = -12.5561
$ws.Range("D39").Value = -7.387900000000001
$ws.Range("C44").Value = -13.21959999999999
$ws.Range("D48").Value = -7.265299999999997
$ws.Range("B51").Value = 5.971900000000004
$ws.Range("C51").Value = -11.92510000000001
$ws.Range("D52").Value = -7.514399999999992
$ws.Range("B55").Value = 4.843799999999996
$ws.Range("D56").Value = -7.907199999999996
$ws.Range("C57").Value = -14.11229999999999
$ws.Range("E57").Value = 16.614
$ws.Range("C63").Value = -11.9474
$ws.Range("A69").Value = -21.6245
$ws.Range("B69").Value = 5.379199999999995
$ws.Range("B70").Value = 6.064600000000008
$ws.Range("C70").Value = -11.577
$ws.Range("A76").Value = -20.03579999999998
$ws.Range("A82").Value = -22.02580000000001
$ws.Range("B83").Value = 5.8346
$ws.Range("D84").Value = -8.923199999999998
$ws.Range("E85").Value = 16.24099999999999
$ws.Range("C99").Value = -12.84219999999999
$ws.Range("D100").Value = -8.356700000000004
$ws.Range("D101").Value = -8.062599999999993
$ws.Range("B102").Value = 8.160000000000007
